$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths.
# Excel's COM ColumnWidth (character units) round-trips into the stored
# OOXML "width" with a constant +5/6 padding offset, so subtract it here
# to land on the exact target stored widths (4, 22, 32).
$ws.Columns.Item(1).ColumnWidth = (4 - 5/6)
$ws.Columns.Item(2).ColumnWidth = (22 - 5/6)
$ws.Columns.Item(3).ColumnWidth = (32 - 5/6)

# Header row (row 2)
$ws.Range("A2").Value = "TT"
$ws.Range("B2").Value = "Ký hiệu chữ viết tắt"
$ws.Range("C2").Value = "Chữ viết đây đủ"

# Data rows - column A values are numeric-looking text, force text via leading apostrophe
$data = @(
    @("1", "TDTT", "Thể dục thể thao"),
    @("2", "HĐND", "Hội đông nhân dân"),
    @("3", "UBND", "ủy ban nhân dân"),
    @("4", "KH", "Kế hoạch"),
    @("5", "ATGT QL 217", "An toàn giao thông Quốc lộ 217"),
    @("6", "UBDS", "Ủy ban dân số"),
    @("7", "THPT", "Trung học phố thông"),
    @("8", "TE", "Trẻ em"),
    @("9", "CĐ-ĐH", "Cao đẳng - đại học"),
    @("10", "CNH- HĐH", "Công nghiệp hóa- hiện đại hóa")
)

$rowIndex = 3
foreach ($rowData in $data) {
    $ws.Range("A$rowIndex").Value = "'" + $rowData[0]
    $ws.Range("B$rowIndex").Value = $rowData[1]
    $ws.Range("C$rowIndex").Value = $rowData[2]
    $rowIndex++
}
